$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    # Force Excel to store the value as text (shared string) even when it
    # looks like a number, without leaving a visible style on the cell.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

# --- Row 5 updates -------------------------------------------------------
# "Fecha Envio" timestamp was re-saved with slightly different precision
$ws.Cells.Item(5, 9).Value2 = 45706.58270331019

# Numero (O5) and Codigo Postal (Q5) are now real numbers instead of text
$ws.Cells.Item(5, 15).Value = 12
$ws.Cells.Item(5, 17).Value = 39001

# Es Serviciable (S5) changed from "Si" to "No"
$ws.Cells.Item(5, 19).Value = "No"

# --- Row 6: a brand-new submission appended at the bottom ---------------
$ws.Cells.Item(6, 2).Value  = "aaa"
Set-TextValue $ws.Cells.Item(6, 3) "666666666"
$ws.Cells.Item(6, 6).Value  = "aaaa"
$ws.Cells.Item(6, 7).Value  = 43.42296284866612
$ws.Cells.Item(6, 8).Value  = -3.683935015

$i6 = $ws.Cells.Item(6, 9)
$i6.NumberFormat = "YYYY-MM-DD HH:MM:SS"
$i6.Value2 = 45706.59697426638

$ws.Cells.Item(6, 10).Value = "N/D"
$ws.Cells.Item(6, 11).Value = "CANTABRIA"
$ws.Cells.Item(6, 12).Value = "RIBAMONTÃN AL MONTE"
$ws.Cells.Item(6, 13).Value = "PONTONES"
$ws.Cells.Item(6, 14).Value = "PONTONES"
Set-TextValue $ws.Cells.Item(6, 15) "99999"
Set-TextValue $ws.Cells.Item(6, 17) "39793"
$ws.Cells.Item(6, 18).Value = "aaaa"
$ws.Cells.Item(6, 19).Value = "No"
$ws.Cells.Item(6, 20).Value = "aaa"
